$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 675.75  # H28
$ws.Cells.Item(28, 9).Value = 693.05554  # I28
$ws.Cells.Item(28, 10).Value = 520  # J28
$ws.Cells.Item(28, 11).Value = 693.05554  # K28
$ws.Cells.Item(28, 12).Value = 520  # L28
$ws.Cells.Item(28, 13).Value = -208.05554  # M28
$ws.Cells.Item(28, 14).Value = -1490  # N28

$ws.Cells.Item(61, 8).Value = 420.2  # H61
$ws.Cells.Item(61, 9).Value = 420.2  # I61
$ws.Cells.Item(61, 11).Value = 1260.6  # K61
$ws.Cells.Item(61, 13).Value = -1088.6  # M61

$ws.Cells.Item(70, 8).Value = 4396.4  # H70
$ws.Cells.Item(70, 9).Value = 2742.6667  # I70
$ws.Cells.Item(70, 10).Value = 6877  # J70
$ws.Cells.Item(70, 11).Value = 8228.000100000001  # K70
$ws.Cells.Item(70, 12).Value = 20631  # L70
$ws.Cells.Item(70, 13).Value = -7958.000100000001  # M70
$ws.Cells.Item(70, 14).Value = -21171  # N70

$ws.Cells.Item(73, 8).Value = 4396.4  # H73
$ws.Cells.Item(73, 9).Value = 2742.6667  # I73
$ws.Cells.Item(73, 10).Value = 6877  # J73
$ws.Cells.Item(73, 11).Value = 8228.000100000001  # K73
$ws.Cells.Item(73, 12).Value = 20631  # L73
$ws.Cells.Item(73, 13).Value = -7292.000100000001  # M73
$ws.Cells.Item(73, 14).Value = -22503  # N73

$ws.Cells.Item(76, 8).Value = 3605.2222  # H76
$ws.Cells.Item(76, 9).Value = 3605.2222  # I76
$ws.Cells.Item(76, 11).Value = 3605.2222  # K76
$ws.Cells.Item(76, 13).Value = -3290.2222  # M76

$ws.Cells.Item(79, 8).Value = 3605.2222  # H79
$ws.Cells.Item(79, 9).Value = 3605.2222  # I79
$ws.Cells.Item(79, 11).Value = 3605.2222  # K79
$ws.Cells.Item(79, 13).Value = -2513.2222  # M79

$ws.Cells.Item(101, 8).Value = 30059628  # H101
$ws.Cells.Item(101, 9).Value = 770206.6  # I101
$ws.Cells.Item(101, 11).Value = 2310619.8  # K101
$ws.Cells.Item(101, 13).Value = -2308997.8  # M101

$ws.Cells.Item(107, 8).Value = 35195.863  # H107
$ws.Cells.Item(107, 9).Value = 538.06665  # I107
$ws.Cells.Item(107, 11).Value = 538.06665  # K107
$ws.Cells.Item(107, 13).Value = 1381.93335  # M107

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3299.459  # H32
$ws.Cells.Item(32, 9).Value = 3172.8167  # I32
$ws.Cells.Item(32, 11).Value = 3172.8167  # K32
$ws.Cells.Item(32, 13).Value = -2885.8167  # M32

$ws.Cells.Item(45, 8).Value = 6905.5415  # H45
$ws.Cells.Item(45, 9).Value = 9167.214  # I45
$ws.Cells.Item(45, 11).Value = 9167.214  # K45
$ws.Cells.Item(45, 13).Value = -8790.214  # M45

$ws.Cells.Item(110, 8).Value = 1499.5555  # H110
$ws.Cells.Item(110, 9).Value = 1592.4286  # I110
$ws.Cells.Item(110, 11).Value = 1592.4286  # K110
$ws.Cells.Item(110, 13).Value = 452.5714  # M110

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 144430.28  # H134
$ws.Cells.Item(134, 9).Value = 177495.53  # I134
$ws.Cells.Item(134, 10).Value = 3903  # J134
$ws.Cells.Item(134, 11).Value = 532486.59  # K134
$ws.Cells.Item(134, 12).Value = 11709  # L134
$ws.Cells.Item(134, 13).Value = -529951.59  # M134
$ws.Cells.Item(134, 14).Value = -16779  # N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 35072.023  # H31
$ws.Cells.Item(31, 9).Value = 31903.295  # I31
$ws.Cells.Item(31, 10).Value = 50463  # J31
$ws.Cells.Item(31, 11).Value = 31903.295  # K31
$ws.Cells.Item(31, 12).Value = 50463  # L31
$ws.Cells.Item(31, 13).Value = -31608.295  # M31
$ws.Cells.Item(31, 14).Value = -51053  # N31

$ws.Cells.Item(34, 8).Value = 35072.023  # H34
$ws.Cells.Item(34, 9).Value = 31903.295  # I34
$ws.Cells.Item(34, 10).Value = 50463  # J34
$ws.Cells.Item(34, 11).Value = 31903.295  # K34
$ws.Cells.Item(34, 12).Value = 50463  # L34
$ws.Cells.Item(34, 13).Value = -31701.295  # M34
$ws.Cells.Item(34, 14).Value = -50867  # N34

$ws.Cells.Item(62, 8).Value = 2570.4285  # H62
$ws.Cells.Item(62, 9).Value = 2873.5  # I62
$ws.Cells.Item(62, 10).Value = 2166.3333  # J62
$ws.Cells.Item(62, 11).Value = 2873.5  # K62
$ws.Cells.Item(62, 12).Value = 2166.3333  # L62
$ws.Cells.Item(62, 13).Value = -2249.5  # M62
$ws.Cells.Item(62, 14).Value = -3414.3333  # N62

$ws.Cells.Item(65, 8).Value = 2570.4285  # H65
$ws.Cells.Item(65, 9).Value = 2873.5  # I65
$ws.Cells.Item(65, 10).Value = 2166.3333  # J65
$ws.Cells.Item(65, 11).Value = 14367.5  # K65
$ws.Cells.Item(65, 12).Value = 10831.6665  # L65
$ws.Cells.Item(65, 13).Value = -11247.5  # M65
$ws.Cells.Item(65, 14).Value = -17071.6665  # N65

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 494.66666  # H5
$ws.Cells.Item(5, 9).Value = 506.5  # I5
$ws.Cells.Item(5, 10).Value = 400  # J5
$ws.Cells.Item(5, 11).Value = 1519.5  # K5
$ws.Cells.Item(5, 12).Value = 1200  # L5
$ws.Cells.Item(5, 13).Value = -1407.5  # M5
$ws.Cells.Item(5, 14).Value = -1424  # N5

$ws.Cells.Item(45, 8).Value = 10000  # H45
$ws.Cells.Item(45, 10).Value = 10000  # J45
$ws.Cells.Item(45, 12).Value = 30000  # L45
$ws.Cells.Item(45, 14).Value = -31064  # N45

$ws.Cells.Item(135, 8).Value = 494.66666  # H135
$ws.Cells.Item(135, 9).Value = 506.5  # I135
$ws.Cells.Item(135, 10).Value = 400  # J135
$ws.Cells.Item(135, 11).Value = 4558.5  # K135
$ws.Cells.Item(135, 12).Value = 3600  # L135
$ws.Cells.Item(135, 13).Value = -2023.5  # M135
$ws.Cells.Item(135, 14).Value = -8670  # N135

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3507.125  # H80
$ws.Cells.Item(80, 9).Value = 2934.5386  # I80
$ws.Cells.Item(80, 11).Value = 2934.5386  # K80
$ws.Cells.Item(80, 13).Value = -1936.5386  # M80

$ws.Cells.Item(83, 8).Value = 3507.125  # H83
$ws.Cells.Item(83, 9).Value = 2934.5386  # I83
$ws.Cells.Item(83, 11).Value = 14672.693  # K83
$ws.Cells.Item(83, 13).Value = -9680.692999999999  # M83

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 300449.84  # H16
$ws.Cells.Item(16, 9).Value = 160539.8  # I16
$ws.Cells.Item(16, 10).Value = 1000000  # J16
$ws.Cells.Item(16, 11).Value = 160539.8  # K16
$ws.Cells.Item(16, 12).Value = 1000000  # L16
$ws.Cells.Item(16, 13).Value = -160369.8  # M16
$ws.Cells.Item(16, 14).Value = -1000340  # N16

$ws.Cells.Item(22, 8).Value = 1060  # H22
$ws.Cells.Item(22, 10).Value = 1100  # J22
$ws.Cells.Item(22, 12).Value = 1100  # L22
$ws.Cells.Item(22, 14).Value = -1690  # N22

$ws.Cells.Item(27, 8).Value = 1060  # H27
$ws.Cells.Item(27, 10).Value = 1100  # J27
$ws.Cells.Item(27, 12).Value = 1100  # L27
$ws.Cells.Item(27, 14).Value = -1314  # N27

$ws.Cells.Item(46, 8).Value = 37326.668  # H46
$ws.Cells.Item(46, 9).Value = 212811.5  # I46
$ws.Cells.Item(46, 11).Value = 212811.5  # K46
$ws.Cells.Item(46, 13).Value = -212623.5  # M46

$ws.Cells.Item(55, 8).Value = 325.61905  # H55
$ws.Cells.Item(55, 10).Value = 350.6  # J55
$ws.Cells.Item(55, 12).Value = 350.6  # L55
$ws.Cells.Item(55, 14).Value = -696.6  # N55

$ws.Cells.Item(68, 8).Value = 0  # H68
$ws.Cells.Item(68, 9).Value = 0  # I68
$ws.Cells.Item(68, 11).Value = 0  # K68
$ws.Cells.Item(68, 13).ClearContents()  # M68

$ws.Cells.Item(71, 8).Value = 0  # H71
$ws.Cells.Item(71, 9).Value = 0  # I71
$ws.Cells.Item(71, 11).Value = 0  # K71
$ws.Cells.Item(71, 13).ClearContents()  # M71

$ws.Cells.Item(82, 8).Value = 1658  # H82
$ws.Cells.Item(82, 9).Value = 1440.3334  # I82
$ws.Cells.Item(82, 10).Value = 2147.75  # J82
$ws.Cells.Item(82, 11).Value = 1440.3334  # K82
$ws.Cells.Item(82, 12).Value = 2147.75  # L82
$ws.Cells.Item(82, 13).Value = -1079.3334  # M82
$ws.Cells.Item(82, 14).Value = -2869.75  # N82

$ws.Cells.Item(85, 8).Value = 1658  # H85
$ws.Cells.Item(85, 9).Value = 1440.3334  # I85
$ws.Cells.Item(85, 10).Value = 2147.75  # J85
$ws.Cells.Item(85, 11).Value = 1440.3334  # K85
$ws.Cells.Item(85, 12).Value = 2147.75  # L85
$ws.Cells.Item(85, 13).Value = -192.3334  # M85
$ws.Cells.Item(85, 14).Value = -4643.75  # N85

$ws.Cells.Item(122, 8).Value = 147217.58  # H122
$ws.Cells.Item(122, 9).Value = 336676  # I122
$ws.Cells.Item(122, 11).Value = 1010028  # K122
$ws.Cells.Item(122, 13).Value = -1007578  # M122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 28999  # H39
$ws.Cells.Item(39, 10).Value = 28999  # J39
$ws.Cells.Item(39, 12).Value = 28999  # L39
$ws.Cells.Item(39, 14).Value = -29825  # N39

$ws.Cells.Item(42, 8).Value = 0  # H42
$ws.Cells.Item(42, 9).Value = 0  # I42
$ws.Cells.Item(42, 11).Value = 0  # K42
$ws.Cells.Item(42, 13).ClearContents()  # M42

$ws.Cells.Item(62, 8).Value = 24322.334  # H62
$ws.Cells.Item(62, 9).Value = 5224.5  # I62
$ws.Cells.Item(62, 11).Value = 5224.5  # K62
$ws.Cells.Item(62, 13).Value = -4600.5  # M62

$ws.Cells.Item(65, 8).Value = 24322.334  # H65
$ws.Cells.Item(65, 9).Value = 5224.5  # I65
$ws.Cells.Item(65, 11).Value = 26122.5  # K65
$ws.Cells.Item(65, 13).Value = -23002.5  # M65

$ws.Cells.Item(126, 8).Value = 2384  # H126
$ws.Cells.Item(126, 9).Value = 2364.6667  # I126
$ws.Cells.Item(126, 10).Value = 2500  # J126
$ws.Cells.Item(126, 11).Value = 7094.000100000001  # K126
$ws.Cells.Item(126, 12).Value = 7500  # L126
$ws.Cells.Item(126, 13).Value = -4624.000100000001  # M126
$ws.Cells.Item(126, 14).Value = -12440  # N126
